$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Changes")

# --- Row 13: merge F13 (Funding) + G13 (Need to get this done. Add more teams)
#     into a single comma-separated string in F13; clear G13; update H13 note.
$ws.Range("F13").Value = "Funding,Need to get this done. Add more teams"
$ws.Range("G13").ClearContents()
$ws.Range("H13").Value = 'Move a card to a new lane and use string after first ","  as wipOverrideComment, if needed'

# --- Row 16: clarify how isBlocked is cleared (blockReason "-" prefix, "" or single char)
$ws.Range("H16").Value = 'Remove isBlocked by starting blockReason with a "-", or by setting it to: "" or single char entry'

# --- Row 20: merge F20 (My GitHub) + G20 (https://github.com/nikantonelli/GroundHog)
#     into a single comma-separated string in F20; clear G20 (and its hyperlink);
#     add explanatory note in H20.
# Removing a hyperlink via the Hyperlinks collection clears ALL hyperlinks on the
# sheet in this runtime, so the other two (F10 -> mailto:nantonelli@planview.com,
# F9 -> mailto:-nantonelli@planview.com) are re-added afterwards.
$ws.Range("F20").Value = "My GitHub, https://github.com/nikantonelli/GroundHog"
$ws.Range("G20").Hyperlinks.Delete()
$ws.Range("G20").ClearContents()
$ws.Range("H20").Value = 'Code looks for last occurence of "," to separate url'

$ws.Hyperlinks.Add($ws.Range("F10"), "mailto:nantonelli@planview.com")
$ws.Hyperlinks.Add($ws.Range("F9"), "mailto:-nantonelli@planview.com")

# --- Re-point view focus to the Changes sheet (activeTab / tabSelected / selection)
$ws.Activate()
$ws.Range("G38").Select()
